$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# Typography sheet: I4 Widget Wildcard Characters = "0-9"
$wsTypo.Range('I4').Value = '0-9'

# Translation sheet: remove the old row 43 (SingleUseId41 / Ok),
# shifting the NFC/MQTT rows up by one
$wsTrans.Rows.Item(43).Delete()

# Append new rows for NFC task interaction / NFC screen / screen timer
$wsTrans.Range('B58').Value = 'SingleUseId61'
$wsTrans.Range('C58').Value = 'PadNumeric'
$wsTrans.Range('D58').Value = 'Left'
$wsTrans.Range('E58').Value = 'LTR'
$wsTrans.Range('F58').Value = '0'

$wsTrans.Range('B59').Value = 'SingleUseId62'
$wsTrans.Range('C59').Value = 'Typographies_button'
$wsTrans.Range('D59').Value = 'Center'
$wsTrans.Range('E59').Value = 'LTR'
$wsTrans.Range('F59').Value = '<value>'

$wsTrans.Range('B60').Value = 'SingleUseId63'
$wsTrans.Range('C60').Value = 'Typographies_button'
$wsTrans.Range('D60').Value = 'Left'
$wsTrans.Range('E60').Value = 'LTR'
$wsTrans.Range('F60').Value = '0'

$wsTrans.Range('B61').Value = 'SingleUseId65'
$wsTrans.Range('C61').Value = 'Typography_label'
$wsTrans.Range('D61').Value = 'Left'
$wsTrans.Range('E61').Value = 'LTR'
$wsTrans.Range('F61').Value = 'Rest waiting time'

$wsTrans.Range('B62').Value = 'SingleUseId66'
$wsTrans.Range('C62').Value = 'Typography_label'
$wsTrans.Range('D62').Value = 'Left'
$wsTrans.Range('E62').Value = 'LTR'
$wsTrans.Range('F62').Value = 'Card read was: '
